$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Whole year")
$ws.Range("C2").Value = 10612736.799
$ws.Range("C3").Value = 9361383.787897468
$ws.Range("C4").Value = 427071.2710208746
$ws.Range("C5").Value = 168391.0278618892
$ws.Range("C6").Value = 1125.222010135986
$ws.Range("C7").Value = 4762481.516371373
$ws.Range("C8").Value = 3134861.941707576
$ws.Range("C9").Value = 37363.77303307541
$ws.Range("C10").Value = 1707840.615235701
$ws.Range("C11").Value = 2916925.832670349
$ws.Range("C12").Value = 9632461.800901605
$ws.Range("C13").Value = 12171623.16534381
$ws.Range("C14").Value = 3576208.440460043
$ws.Range("C15").Value = 48166.72185764778
$ws.Range("C16").Value = 8191767.866063803
$ws.Range("C17").Value = 390241.5577952634
$ws.Range("C18").Value = 15918445.0859807
$ws.Range("C19").Value = 25570.9059762565
$ws.Range("C20").Value = 8063011.671290446
$ws.Range("C21").Value = 1076708.486683669
$ws.Range("C22").Value = 9222530.238275034
$ws.Range("C23").Value = 16911643.5345532
$ws.Range("C24").Value = 10056633.80555829
$ws.Range("C25").Value = 4864814.490741099

$ws = $wb.Worksheets.Item("Winter")
$ws.Range("C2").Value = 2632858.689
$ws.Range("C3").Value = 3503815.793273349
$ws.Range("C4").Value = 84927.53610156653
$ws.Range("C5").Value = 90142.40046959631
$ws.Range("C6").Value = 430.7349768272641
$ws.Range("C7").Value = 2267835.727646634
$ws.Range("C8").Value = 422996.1914057396
$ws.Range("C9").Value = 4906.550578521911
$ws.Range("C10").Value = 1202965.617614802
$ws.Range("C11").Value = 1633985.851825688
$ws.Range("C12").Value = 70729.94895657552
$ws.Range("C13").Value = 83529.94538636699
$ws.Range("C14").Value = 1993490.370876481
$ws.Range("C15").Value = 134.0736238995695
$ws.Range("C16").Value = 4574555.5783332
$ws.Range("C17").Value = 47037.37409917702
$ws.Range("C18").Value = 861371.2632042464
$ws.Range("C19").Value = 25570.6814868186
$ws.Range("C20").Value = 84848.62563179298
$ws.Range("C21").Value = 470922.3242761792
$ws.Range("C22").Value = 6832519.15028209
$ws.Range("C23").Value = 846992.9392212338
$ws.Range("C24").Value = 7167326.040527486
$ws.Range("C25").Value = 2030911.169037429

$ws = $wb.Worksheets.Item("Spring")
$ws.Range("C2").Value = 3611313.025
$ws.Range("C3").Value = 2503991.560105643
$ws.Range("C4").Value = 137509.8244841086
$ws.Range("C5").Value = 36926.52602441507
$ws.Range("C6").Value = 303.7396682196943
$ws.Range("C7").Value = 1013057.410638321
$ws.Range("C8").Value = 895009.8610308161
$ws.Range("C9").Value = 17362.58003323062
$ws.Range("C10").Value = 300190.6254883407
$ws.Range("C11").Value = 555922.9406147172
$ws.Range("C12").Value = 3427843.326436596
$ws.Range("C13").Value = 4276221.634823779
$ws.Range("C14").Value = 657593.8075766476
$ws.Range("C16").Value = 3007150.016121518
$ws.Range("C17").Value = 120028.7507701389
$ws.Range("C18").Value = 5095320.870198407
$ws.Range("C20").Value = 3614326.961078656
$ws.Range("C21").Value = 392742.8991128182
$ws.Range("C22").Value = 1350821.221885394
$ws.Range("C23").Value = 5554481.332280847
$ws.Range("C24").Value = 1447241.690690052
$ws.Range("C25").Value = 1383142.843993499

$ws = $wb.Worksheets.Item("Summer")
$ws.Range("C2").Value = 2985749.657
$ws.Range("C3").Value = 1936567.557996062
$ws.Range("C4").Value = 132145.0986409418
$ws.Range("C5").Value = 18710.80647416027
$ws.Range("C6").Value = 217.8549673544717
$ws.Range("C7").Value = 946717.9214573825
$ws.Range("C8").Value = 1431970.464481312
$ws.Range("C9").Value = 9574.845410824915
$ws.Range("C10").Value = 68165.67793764456
$ws.Range("C11").Value = 0
$ws.Range("C12").Value = 4876128.929057013
$ws.Range("C13").Value = 6374173.245414898
$ws.Range("C14").Value = 0
$ws.Range("C15").Value = 48032.64823374821
$ws.Range("C16").Value = 177464.7296604221
$ws.Range("C17").Value = 102803.3748363527
$ws.Range("C18").Value = 7576664.279569159
$ws.Range("C19").Value = 0.2244894378993121
$ws.Range("C20").Value = 2839789.764353736
$ws.Range("C21").Value = 89883.38501352376
$ws.Range("C22").Value = 113923.2334113341
$ws.Range("C23").Value = 7982391.865773327
$ws.Range("C24").Value = 461273.8499141681
$ws.Range("C25").Value = 803194.6726352386

$ws = $wb.Worksheets.Item("Fall")
$ws.Range("C2").Value = 1382815.428
$ws.Range("C3").Value = 1417008.876522413
$ws.Range("C4").Value = 72488.81179425762
$ws.Range("C5").Value = 22611.29489371755
$ws.Range("C6").Value = 172.8923977345558
$ws.Range("C7").Value = 534870.456629036
$ws.Range("C8").Value = 384885.4247897079
$ws.Range("C9").Value = 5519.797010497965
$ws.Range("C10").Value = 136518.6941949135
$ws.Range("C11").Value = 726971.4753909551
$ws.Range("C12").Value = 1257730.313854475
$ws.Range("C13").Value = 1437667.017306131
$ws.Range("C14").Value = 925067.9555927777
$ws.Range("C16").Value = 432531.2941499744
$ws.Range("C17").Value = 120359.2339824744
$ws.Range("C18").Value = 2385057.560746178
$ws.Range("C20").Value = 1524032.102866277
$ws.Range("C21").Value = 123150.858414662
$ws.Range("C22").Value = 925174.3813804917
$ws.Range("C23").Value = 2527777.39727779
$ws.Range("C24").Value = 980792.2244265833
$ws.Range("C25").Value = 647565.8050749324

Write-Host "Updated sankey values for all sheets"